$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $oldStyle = $Range.Style
    $Range.Value = "'" + $Text
    $Range.Style = $oldStyle
}

$ws.Range("D2").Value = "69.102.68"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "3.769.86"
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "632.56"
$ws.Range("E5").Value = "  +3.32%  "
Set-TextValue $ws.Range("D6") "166.37"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "3.768.34"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.43%  "
Set-TextValue $ws.Range("D10") "0.158"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +1.51%  "
Set-TextValue $ws.Range("D12") "6.79"
$ws.Range("E12").Value = "  +0.28%  "
Set-TextValue $ws.Range("D13") "0.0000240"
$ws.Range("E13").Value = "  -3.20%  "
Set-TextValue $ws.Range("D14") "34.88"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "4.406.15"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "3.769.10"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "69.126.83"
$ws.Range("E17").Value = "  +0.76%  "
Set-TextValue $ws.Range("D18") "17.63"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("E19").Value = "  +0.14%  "
Set-TextValue $ws.Range("D20") "7.02"
$ws.Range("E20").Value = "  -1.67%  "
Set-TextValue $ws.Range("D21") "461.99"
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  -1.61%  "
Set-TextValue $ws.Range("D23") "0.707"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  -2.42%  "
Set-TextValue $ws.Range("D25") "82.09"
$ws.Range("E25").Value = "  -1.96%  "
Set-TextValue $ws.Range("D26") "12.11"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +0.55%  "
Set-TextValue $ws.Range("D28") "10.08"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D30").Value = "3.920.73"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  +2.90%  "
Set-TextValue $ws.Range("D32") "2.69"
$ws.Range("E32").Value = "  +2.04%  "
Set-TextValue $ws.Range("D33") "7.06"
$ws.Range("E33").Value = "  -2.98%  "
Set-TextValue $ws.Range("D34") "0.177"
$ws.Range("E34").Value = "  +20.44%  "
Set-TextValue $ws.Range("D35") "28.43"
$ws.Range("E35").Value = "  -1.98%  "
Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").Value = "3.724.28"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  +4.51%  "
$ws.Range("E41").Value = "  -1.81%  "
Set-TextValue $ws.Range("D42") "0.963"
$ws.Range("E42").Value = "  -1.87%  "
Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +0.08%  "
Set-TextValue $ws.Range("D45") "157.47"
$ws.Range("E45").Value = "  +2.21%  "
Set-TextValue $ws.Range("D46") "1.97"
$ws.Range("E46").Value = "  +5.24%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D48") "47.04"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D49") "43.02"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -0.28%  "
